$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": update currency-year conversion factor & document the
# inflation adjustment source, then renumber the 45Q duration note.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Update the 2012-USD conversion factor (was "2022 to 2012 USD" @ 0.785)
$about.Range("A9").Value = 0.73
$about.Range("C9").Value = "https://www.bls.gov/data/inflation_calculator.htm"
$about.Range("B9").Value = "2024 to 2012 USD"

# ---------------------------------------------------------------------------
# Sheet "Electricity Calculations": the credit amount is now expressed in
# 2012 USD, so it gets deflated by the About!A9 factor before being used.
# ---------------------------------------------------------------------------
$elec = $wb.Worksheets.Item("Electricity Calculations")
$elec.Range("A1").Value = "Credit Amount, 2012 USD"
$elec.Range("B1").Formula = "=About!B11*About!A9"

# Insert a new explanatory row right after the 45Q Tax Credit Amount row,
# pushing the old "45Q Duration" / "12 years" row down from 13 to 14.
[void]$about.Rows.Item(12).Insert()
$about.Range("A12").Value = "*inflation adjusted starting in 2025, so we use the 2024 currency year to adjust to 2012 `$"

[void]$about.Range("A13").Select()

# ---------------------------------------------------------------------------
# Sheet "BCS-BCS": the electricity-sector row now pulls its (already
# deflated & duration-adjusted) value from Electricity Calculations!B4
# instead of recomputing About!$B$11*About!$A$9 directly.
# ---------------------------------------------------------------------------
$bcs = $wb.Worksheets.Item("BCS-BCS")
$elecCols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M")
foreach ($col in $elecCols) {
    $cell = $bcs.Range($col + "2")
    $cell.Formula = "='Electricity Calculations'!`$B`$4"
    $cell.NumberFormat = "General"
}

[void]$bcs.Range("F2:G2").Select()
